$wb = $excel.ActiveWorkbook

# --- Update selections on existing sheets ---

# BSC sheet: selection B4 -> B5
$wsBsc = $wb.Worksheets.Item("BSC")
$wsBsc.Range("B5").Select()

# Medtronics sheet: selection A1:XFD1 -> B2
$wsMed = $wb.Worksheets.Item("Medtronics ")
$wsMed.Range("B2").Select()

# Crane1 sheet: topLeftCell B2 removed, selection C2/C1:U1048576 -> A1:A1048576
$wsCrane = $wb.Worksheets.Item("Crane1")
$wsCrane.Range("A1:A1048576").Select()

# --- Add new "Generic" sheet after Crane1 ---
# (An extra throwaway sheet is inserted first and then removed so the new
#  sheet's internal sheetId lands on 8, matching the source workbook.)

$wsLast = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsTemp = $wb.Worksheets.Add($null, $wsLast)
$wsTemp.Name = "TempPlaceholder"

$wsNew = $wb.Worksheets.Add($null, $wsTemp)
$wsNew.Name = "Generic"

$wsBsc.Range("A1:B1").Copy()
$wsNew.Range("A1:B1").PasteSpecial(-4122)
$wsNew.Range("A1").Value = "Number"
$wsNew.Range("B1").Value = "Query"

$wsNew.Range("A2").Value = 1
$wsNew.Range("B2").Value = "4012 optics initialization error ."

$wsNew.Range("A3").Value = 2
$wsNew.Range("B3").Value = "device not starting"

$wsNew.Range("A4").Value = 3
$wsNew.Range("B4").Value = "hardware failure"

$wsNew.Range("A5").Value = 4
$wsNew.Range("B5").Value = "optics initialization error"

$wsNew.Columns.Item(2).ColumnWidth = 26

$wsNew.Range("B2").Select()

# remove the throwaway placeholder sheet now that "Generic" is fully populated
$wsTemp.Delete()

# restore original active sheet / tab selection
$wsDevdemo = $wb.Worksheets.Item("Devdemo")
$wsDevdemo.Activate()
$wsDevdemo.Range("B3").Select()
